$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "questions per category" value in B7
$ws.Range("B7").Value = "5, 5, 5, 5, 3, 3, 3"

# Update the selection/active cell to B7 as shown in the saved view
$ws.Range("B7").Select()
